$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) column at C - the wide spacer column that used to be C
# (with its 22.8-wide custom width) becomes D; the value columns that used to
# be D (Possible) and E (Score) shift right to E and F. Column A (category
# labels) and B (unused spacer) are untouched.
$ws.Columns("C:C").Insert()

# Insert 2 new rows right after the "Web app 1 functionality" row (old row 9)
# to make room for the "Array" and "Functions" sub-rows.
$ws.Rows("10:11").Insert()

# Insert 3 new rows right after the "Web app 2 functionality" row (now row 12)
# to make room for the "Arrays", "Functions", and "Style and best practices" sub-rows.
$ws.Rows("13:15").Insert()

# Move the blank bold-styled cell from A3 to B3
$ws.Range("A3").Clear()
$ws.Range("B3").Font.Bold = $true

# --- Fill in the new label text (this order matches the order the labels
#     were originally typed in, as reflected by the shared-strings table) ---
$ws.Range("A15").Value = "Style and best practices"
$ws.Range("B10").Value = "Array"
$ws.Range("B11").Value = "Functions"
$ws.Range("B13").Value = "Arrays"
$ws.Range("B14").Value = "Functions"
$ws.Range("B22").Value = "Grade book, Price list, To-do list"
$ws.Range("B20").Value = "Roman to decimal converter, Average of scores, Decimal to roman numeral converter"

# --- Fill in the new sub-rows' numeric values under "Web app 1 functionality" ---
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 2

$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 6

# --- Update "Web app 1 functionality" possible/score values (row 9) ---
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 4

# --- Update "Web app 2 functionality" possible/score values (row 12) ---
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 7

# --- Fill in the new sub-rows' numeric values under "Web app 2 functionality" ---
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 4

$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 8

$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 5

# --- Add the new "used in" notes below the table ---
$ws.Range("A20").Value = 1
$ws.Range("A22").Value = 2

# Narrow decorative spacer columns at the left (A) and right (G) edges of the table
$ws.Columns("A:A").ColumnWidth = 1.8
$ws.Columns("G:G").ColumnWidth = 0.5

$ws.Range("H14").Select()
